$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells hold text-formatted numbers (e.g. "175.10", "70.620.51") that
# must stay literal strings instead of being auto-converted to numeric values by
# Excel. Temporarily switch each cell to text format before assigning the new
# value, then clear the (now redundant) explicit formatting so the cell keeps
# its original default style, matching the source workbook's layout.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.564.90'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.622.21'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.10'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.631'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.612.31'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.193'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.61'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +13.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.615'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '48.36'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.41%  '
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '683.17'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.213.35'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '9.00'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.627.56'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.620.24'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("E20").Value = '  -0.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.74'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.42'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.934'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.05'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.66'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.90'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.95%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.75'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.36'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.12'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.43'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.37'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.95'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '570.65'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.07'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.107'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '58.48'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0451'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.59%  '
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.346'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.532.07'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.62%  '
$ws.Range("E44").Value = '  -2.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.19'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.07%  '
$ws.Range("E46").Value = '  -5.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.99'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +5.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.65'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.59%  '
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.85'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.40%  '
$ws.Range("E51").Value = '  -4.09%  '
